## Applies the "driftmove2.pptx" commit:
##  1. Refresh the cached "datetimeFigureOut" header/footer date field
##     (2/20/17 -> 2/24/17) on the slide master and every slide layout.
##  2. Add a new fully-opaque "Oval 29" circle (duplicated from the
##     existing "Oval 97" circle, recolored/repositioned) at the front
##     of the slide's z-order.
##  3. Remove the old semi-transparent "Oval 97" circle it replaces.
##  4. Nudge "TextBox 111" up a bit.
##  5. Recolor "Rectangle 83" and give it an explicit no-line outline.

function Find-ShapeByName {
    param($Shapes, [string]$Name)
    for ($i = 1; $i -le $Shapes.Count; $i++) {
        $shp = $Shapes.Item($i)
        if ($shp.Name -eq $Name) {
            return $shp
        }
    }
    return $null
}

# EMU -> points, nudged by half an EMU so the float32 round-trip inside
# the host lands on the exact target EMU value instead of one EMU short.
function Emu2Pt {
    param([double]$Emu)
    return ($Emu + 0.5) / 12700.0
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Date placeholder ("datetimeFigureOut" field) on master + layouts.
# ---------------------------------------------------------------------
$design = $p.Designs.Item(1)
$master = $design.SlideMaster

function Update-DatePlaceholder {
    param($Shapes)
    for ($i = 1; $i -le $Shapes.Count; $i++) {
        $shp = $Shapes.Item($i)
        if ($shp.HasTextFrame) {
            try {
                if ($shp.PlaceholderFormat.Type -eq 16) {
                    if ($shp.TextFrame.TextRange.Text -eq "2/20/17") {
                        $shp.TextFrame.TextRange.Text = "2/24/17"
                    }
                }
            } catch {
            }
        }
    }
}

Update-DatePlaceholder $master.Shapes
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $cl = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $cl.Shapes
}

# ---------------------------------------------------------------------
# 2 & 3. Swap "Oval 97" (id 98) for a new, fully opaque "Oval 29".
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

$oldOval = Find-ShapeByName $s.Shapes "Oval 97"

$dupRange = $oldOval.Duplicate()
$newOval = $dupRange.Item(1)
$newOval.Name = "Oval 29"

$newOval.Left = Emu2Pt 8496787
$newOval.Top = Emu2Pt 1196074
$newOval.Width = Emu2Pt 547910
$newOval.Height = Emu2Pt 548640

$newOval.Fill.ForeColor.ObjectThemeColor = 5   # msoThemeColorAccent1
$newOval.Fill.Transparency = 0
$newOval.Line.Visible = $false

$newOval.ZOrder(1)   # msoSendToBack -> becomes the first shape in the tree

$oldOval.Delete()

# ---------------------------------------------------------------------
# 4. Move "TextBox 111" (id 112) up slightly.
# ---------------------------------------------------------------------
$tb111 = Find-ShapeByName $s.Shapes "TextBox 111"
$tb111.Top = Emu2Pt 1123255

# ---------------------------------------------------------------------
# 5. Recolor "Rectangle 83" (id 84) and force an explicit no-line.
# ---------------------------------------------------------------------
$rect83 = Find-ShapeByName $s.Shapes "Rectangle 83"
$rect83.Fill.ForeColor.RGB = 0x232563   # BGR for srgbClr val="632523"
$rect83.Line.Visible = $false
